# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates the "K" column (column G) values on the active sheet for the
# rows whose strikeout counts changed when the save data was regenerated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 2
    4  = 1
    5  = 2
    6  = 0
    8  = 0
    9  = 1
    10 = 1
    11 = 0
    12 = 2
    13 = 0
    14 = 1
    15 = 4
    16 = 3
    17 = 0
    18 = 0
    19 = 2
    20 = 0
    21 = 1
    22 = 0
    23 = 4
    24 = 0
    25 = 2
    26 = 0
    27 = 0
    28 = 2
    29 = 0
    30 = 0
    31 = 2
    32 = 0
    33 = 1
    34 = 3
    35 = 3
    36 = 0
    37 = 1
    39 = 1
    40 = 0
    43 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
